# Insert a new weekly price record (row 251) into the Mango price history
# sheet. All existing rows from 251 downward shift down by one row; the new
# row carries a fresh observation (date serial 44988 = 2023-03-03) with
# slightly higher prices than the following (previously-251) row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 251, shifting rows 251:334 down to 252:335.
$ws.Rows("251:251").Insert()

# Populate the newly inserted row 251 with the new observation.
$ws.Range("A251").Value = 4
$ws.Range("B251").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C251").Value = "Los Lagos"
$ws.Range("D251").Value = 44988
$ws.Range("E251").Value = 10
$ws.Range("F251").Value = "Fruta"
$ws.Range("G251").Value = 100108
$ws.Range("H251").Value = "Tropicales y subtropicales"
$ws.Range("I251").Value = 100108002
$ws.Range("J251").Value = "Mango"
$ws.Range("K251").Value = "Sin especificar"
$ws.Range("L251").Value = "Primera"
$ws.Range("M251").Value = 200
$ws.Range("N251").Value = 8500
$ws.Range("O251").Value = 9000
$ws.Range("P251").Value = 8750
$ws.Range("Q251").Value = "`$/bandeja 4 kilos"
$ws.Range("R251").Value = "Perú"
$ws.Range("S251").Value = 2188
$ws.Range("T251").Value = 4
